# Scheduled-runner style market-price refresh for the Famfrit profit sheets.
# Updates cached currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# per-row across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4,8).Value = 181.66667
$ws.Cells.Item(33,8).Value = 12330.277
$ws.Cells.Item(33,9).Value = 14374.934
$ws.Cells.Item(33,11).Value = 14374.934
$ws.Cells.Item(33,13).Value = -14145.934
$ws.Cells.Item(43,8).Value = 1080.5555
$ws.Cells.Item(43,9).Value = 925.5
$ws.Cells.Item(43,10).Value = 1274.375
$ws.Cells.Item(43,11).Value = 925.5
$ws.Cells.Item(43,12).Value = 1274.375
$ws.Cells.Item(43,13).Value = -856.5
$ws.Cells.Item(43,14).Value = -1412.375
$ws.Cells.Item(116,8).Value = 3322.2
$ws.Cells.Item(116,9).Value = 3002.5
$ws.Cells.Item(116,10).Value = 3535.3333
$ws.Cells.Item(116,11).Value = 3002.5
$ws.Cells.Item(116,12).Value = 3535.3333
$ws.Cells.Item(116,13).Value = 439.5
$ws.Cells.Item(116,14).Value = -10419.3333
$ws.Cells.Item(138,8).Value = 6784.2827
$ws.Cells.Item(138,10).Value = 6933.5684
$ws.Cells.Item(138,12).Value = 20800.7052
$ws.Cells.Item(138,14).Value = -31080.7052

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(41,8).Value = 8044.467
$ws.Cells.Item(41,9).Value = 4320.5386
$ws.Cells.Item(41,10).Value = 32250
$ws.Cells.Item(41,11).Value = 4320.5386
$ws.Cells.Item(41,12).Value = 32250
$ws.Cells.Item(41,13).Value = -3906.5386
$ws.Cells.Item(41,14).Value = -33078
$ws.Cells.Item(74,8).Value = 20432286
$ws.Cells.Item(74,9).Value = 24418644
$ws.Cells.Item(74,10).Value = 2208.5
$ws.Cells.Item(74,11).Value = 24418644
$ws.Cells.Item(74,12).Value = 2208.5
$ws.Cells.Item(74,13).Value = -24417770
$ws.Cells.Item(74,14).Value = -3956.5
$ws.Cells.Item(77,8).Value = 20432286
$ws.Cells.Item(77,9).Value = 24418644
$ws.Cells.Item(77,10).Value = 2208.5
$ws.Cells.Item(77,11).Value = 122093220
$ws.Cells.Item(77,12).Value = 11042.5
$ws.Cells.Item(77,13).Value = -122088852
$ws.Cells.Item(77,14).Value = -19778.5
$ws.Cells.Item(97,8).Value = 1206.931
$ws.Cells.Item(97,9).Value = 1288.5769
$ws.Cells.Item(97,10).Value = 499.33334
$ws.Cells.Item(97,11).Value = 1288.5769
$ws.Cells.Item(97,12).Value = 499.33334
$ws.Cells.Item(97,13).Value = -792.5769
$ws.Cells.Item(97,14).Value = -1491.33334
$ws.Cells.Item(119,8).Value = 100000
$ws.Cells.Item(119,10).Value = 100000
$ws.Cells.Item(119,12).Value = 100000
$ws.Cells.Item(119,14).Value = -109676
$ws.Cells.Item(132,8).Value = 25680616
$ws.Cells.Item(132,9).Value = 4167.067
$ws.Cells.Item(132,11).Value = 12501.201
$ws.Cells.Item(132,13).Value = -9971.201000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20,8).Value = 1176.9736
$ws.Cells.Item(20,9).Value = 892.38464
$ws.Cells.Item(20,11).Value = 892.38464
$ws.Cells.Item(20,13).Value = -645.38464
$ws.Cells.Item(42,8).Value = 200684
$ws.Cells.Item(42,10).Value = 200684
$ws.Cells.Item(42,12).Value = 200684
$ws.Cells.Item(42,14).Value = -201340
$ws.Cells.Item(69,8).Value = 100000
$ws.Cells.Item(69,9).Value = 100000
$ws.Cells.Item(69,11).Value = 100000
$ws.Cells.Item(69,13).Value = -99189
$ws.Cells.Item(72,8).Value = 100000
$ws.Cells.Item(72,9).Value = 100000
$ws.Cells.Item(72,11).Value = 300000
$ws.Cells.Item(72,13).Value = -295944
$ws.Cells.Item(81,8).Value = 21530.5
$ws.Cells.Item(81,10).Value = 21530.5
$ws.Cells.Item(81,12).Value = 21530.5
$ws.Cells.Item(81,14).Value = -23652.5
$ws.Cells.Item(84,8).Value = 21530.5
$ws.Cells.Item(84,10).Value = 21530.5
$ws.Cells.Item(84,12).Value = 64591.5
$ws.Cells.Item(84,14).Value = -75199.5
$ws.Cells.Item(86,8).Value = 5608.4644
$ws.Cells.Item(86,9).Value = 6424.773
$ws.Cells.Item(86,10).Value = 2615.3333
$ws.Cells.Item(86,11).Value = 6424.773
$ws.Cells.Item(86,12).Value = 2615.3333
$ws.Cells.Item(86,13).Value = -5301.773
$ws.Cells.Item(86,14).Value = -4861.3333
$ws.Cells.Item(89,8).Value = 5608.4644
$ws.Cells.Item(89,9).Value = 6424.773
$ws.Cells.Item(89,10).Value = 2615.3333
$ws.Cells.Item(89,11).Value = 32123.865
$ws.Cells.Item(89,12).Value = 13076.6665
$ws.Cells.Item(89,13).Value = -26507.865
$ws.Cells.Item(89,14).Value = -24308.6665
$ws.Cells.Item(99,8).Value = 2592.75
$ws.Cells.Item(99,9).Value = 1549
$ws.Cells.Item(99,10).Value = 4332.3335
$ws.Cells.Item(99,11).Value = 1549
$ws.Cells.Item(99,12).Value = 4332.3335
$ws.Cells.Item(99,13).Value = -51
$ws.Cells.Item(99,14).Value = -7328.3335
$ws.Cells.Item(105,8).Value = 8701.412
$ws.Cells.Item(105,9).Value = 18031.285
$ws.Cells.Item(105,11).Value = 18031.285
$ws.Cells.Item(105,13).Value = -16284.285
$ws.Cells.Item(109,8).Value = 84999.5
$ws.Cells.Item(109,9).Value = 49999
$ws.Cells.Item(109,10).Value = 120000
$ws.Cells.Item(109,11).Value = 49999
$ws.Cells.Item(109,12).Value = 120000
$ws.Cells.Item(109,13).Value = -48612
$ws.Cells.Item(109,14).Value = -122774
$ws.Cells.Item(141,8).Value = 119950
$ws.Cells.Item(141,10).Value = 119950
$ws.Cells.Item(141,12).Value = 119950
$ws.Cells.Item(141,14).Value = -130310

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 3022.4285
$ws.Cells.Item(31,9).Value = 1969.875
$ws.Cells.Item(31,10).Value = 3811.8438
$ws.Cells.Item(31,11).Value = 1969.875
$ws.Cells.Item(31,12).Value = 3811.8438
$ws.Cells.Item(31,13).Value = -1674.875
$ws.Cells.Item(31,14).Value = -4401.843800000001
$ws.Cells.Item(34,8).Value = 3022.4285
$ws.Cells.Item(34,9).Value = 1969.875
$ws.Cells.Item(34,10).Value = 3811.8438
$ws.Cells.Item(34,11).Value = 1969.875
$ws.Cells.Item(34,12).Value = 3811.8438
$ws.Cells.Item(34,13).Value = -1767.875
$ws.Cells.Item(34,14).Value = -4215.843800000001
$ws.Cells.Item(107,8).Value = 823.3333
$ws.Cells.Item(107,9).Value = 1051.75
$ws.Cells.Item(107,10).Value = 640.6
$ws.Cells.Item(107,11).Value = 1051.75
$ws.Cells.Item(107,12).Value = 640.6
$ws.Cells.Item(107,13).Value = 868.25
$ws.Cells.Item(107,14).Value = -4480.6
$ws.Cells.Item(131,8).Value = 34999.5
$ws.Cells.Item(131,10).Value = 34999.5
$ws.Cells.Item(131,12).Value = 34999.5
$ws.Cells.Item(131,14).Value = -45079.5
$ws.Cells.Item(134,8).Value = 3038.3096
$ws.Cells.Item(134,9).Value = 3084.375
$ws.Cells.Item(134,11).Value = 9253.125
$ws.Cells.Item(134,13).Value = -6718.125
$ws.Cells.Item(141,8).Value = 110787.125
$ws.Cells.Item(141,10).Value = 122471
$ws.Cells.Item(141,12).Value = 122471
$ws.Cells.Item(141,14).Value = -132831

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40,8).Value = 371.75
$ws.Cells.Item(40,9).Value = 45.166668
$ws.Cells.Item(40,10).Value = 698.3333
$ws.Cells.Item(40,11).Value = 180.666672
$ws.Cells.Item(40,12).Value = 2793.3332
$ws.Cells.Item(40,13).Value = -111.666672
$ws.Cells.Item(40,14).Value = -2931.3332
$ws.Cells.Item(121,8).Value = 475
$ws.Cells.Item(121,10).Value = 0
$ws.Cells.Item(121,12).Value = 0
$ws.Cells.Item(121,14).ClearContents()
$ws.Cells.Item(128,8).Value = 156015
$ws.Cells.Item(128,9).Value = 156015
$ws.Cells.Item(128,11).Value = 468045
$ws.Cells.Item(128,13).Value = -463065
$ws.Cells.Item(131,8).Value = 41481.484
$ws.Cells.Item(131,10).Value = 10052
$ws.Cells.Item(131,12).Value = 30156
$ws.Cells.Item(131,14).Value = -40236
$ws.Cells.Item(134,8).Value = 2109.35
$ws.Cells.Item(134,10).Value = 8509
$ws.Cells.Item(134,12).Value = 25527
$ws.Cells.Item(134,14).Value = -35667
$ws.Cells.Item(137,8).Value = 2099.1924
$ws.Cells.Item(137,9).Value = 1087.1666
$ws.Cells.Item(137,10).Value = 2966.6428
$ws.Cells.Item(137,11).Value = 3261.4998
$ws.Cells.Item(137,12).Value = 8899.928400000001
$ws.Cells.Item(137,13).Value = 1838.5002
$ws.Cells.Item(137,14).Value = -19099.9284

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107,8).Value = 998.55554
$ws.Cells.Item(107,10).Value = 998
$ws.Cells.Item(107,12).Value = 998
$ws.Cells.Item(107,14).Value = -4838
$ws.Cells.Item(119,8).Value = 61069
$ws.Cells.Item(119,10).Value = 61069
$ws.Cells.Item(119,12).Value = 61069
$ws.Cells.Item(119,14).Value = -70745
$ws.Cells.Item(122,8).Value = 4373.125
$ws.Cells.Item(122,9).Value = 3696.6667
$ws.Cells.Item(122,10).Value = 4779
$ws.Cells.Item(122,11).Value = 11090.0001
$ws.Cells.Item(122,12).Value = 14337
$ws.Cells.Item(122,13).Value = -8640.000100000001
$ws.Cells.Item(122,14).Value = -19237
$ws.Cells.Item(139,8).Value = 99987.5
$ws.Cells.Item(139,10).Value = 99987.5
$ws.Cells.Item(139,12).Value = 99987.5
$ws.Cells.Item(139,14).Value = -110267.5
$ws.Cells.Item(6,8).Value = 39995
$ws.Cells.Item(6,10).Value = 39995
$ws.Cells.Item(6,12).Value = 39995
$ws.Cells.Item(6,14).Value = -40219

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 7271.773
$ws.Cells.Item(7,9).Value = 5099.4
$ws.Cells.Item(7,10).Value = 7910.706
$ws.Cells.Item(7,11).Value = 5099.4
$ws.Cells.Item(7,12).Value = 7910.706
$ws.Cells.Item(7,13).Value = -4987.4
$ws.Cells.Item(7,14).Value = -8134.706
$ws.Cells.Item(17,8).Value = 8030.625
$ws.Cells.Item(17,9).Value = 7356.4287
$ws.Cells.Item(17,11).Value = 7356.4287
$ws.Cells.Item(17,13).Value = -7186.4287
$ws.Cells.Item(40,8).Value = 7326
$ws.Cells.Item(40,9).Value = 5989.5
$ws.Cells.Item(40,10).Value = 9999
$ws.Cells.Item(40,11).Value = 5989.5
$ws.Cells.Item(40,12).Value = 9999
$ws.Cells.Item(40,13).Value = -5853.5
$ws.Cells.Item(40,14).Value = -10271
$ws.Cells.Item(46,8).Value = 2016.8718
$ws.Cells.Item(46,9).Value = 952.0714
$ws.Cells.Item(46,10).Value = 4727.273
$ws.Cells.Item(46,11).Value = 952.0714
$ws.Cells.Item(46,12).Value = 4727.273
$ws.Cells.Item(46,13).Value = -764.0714
$ws.Cells.Item(46,14).Value = -5103.273
$ws.Cells.Item(93,8).Value = 1867.2174
$ws.Cells.Item(93,9).Value = 1467.9143
$ws.Cells.Item(93,10).Value = 3137.7273
$ws.Cells.Item(93,11).Value = 1467.9143
$ws.Cells.Item(93,12).Value = 3137.7273
$ws.Cells.Item(93,13).Value = -219.9142999999999
$ws.Cells.Item(93,14).Value = -5633.7273
$ws.Cells.Item(126,8).Value = 7271.773
$ws.Cells.Item(126,9).Value = 5099.4
$ws.Cells.Item(126,10).Value = 7910.706
$ws.Cells.Item(126,11).Value = 15298.2
$ws.Cells.Item(126,12).Value = 23732.118
$ws.Cells.Item(126,13).Value = -12828.2
$ws.Cells.Item(126,14).Value = -28672.118
$ws.Cells.Item(132,8).Value = 4313.2354
$ws.Cells.Item(132,9).Value = 4077.5833
$ws.Cells.Item(132,11).Value = 12232.7499
$ws.Cells.Item(132,13).Value = -9702.749899999999
$ws.Cells.Item(134,8).Value = 100000
$ws.Cells.Item(134,10).Value = 100000
$ws.Cells.Item(134,12).Value = 100000
$ws.Cells.Item(134,14).Value = -110140

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96,8).Value = 5820
$ws.Cells.Item(96,9).Value = 5150.125
$ws.Cells.Item(96,11).Value = 5150.125
$ws.Cells.Item(96,13).Value = -3777.125
